$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New customer info rows 9-13 (entered first, so they claim the lower shared-string indices)
$ws.Range("B9").Value = "kundennummer "
$ws.Range("B10").Value = "vorname"
$ws.Range("B11").Value = "nachname"
$ws.Range("B12").Value = "straße hausnummer"
$ws.Range("B13").Value = "plz ort"

# Row 7: new header row, split into separate fields with a bottom border style
$ws.Range("C7").Value = "straße"
$ws.Range("E7").Value = "plz"
$ws.Range("F7").Value = "ort"
$ws.Range("H7").Value = "ort,"
$ws.Range("B7").Value = "firmenname,"
$ws.Range("D7").Value = "hausnummer,"
$ws.Range("I7").Value = " den datum"

# Apply bottom-border style to the B7:F7 header cells
$ws.Range("B7:F7").Borders.Item(9).LineStyle = 1

# Apply bottom-border style to B23 and B25 (zusatzkosten / mwst prozent rows)
$ws.Range("B23").Borders.Item(9).LineStyle = 1
$ws.Range("B25").Borders.Item(9).LineStyle = 1

# Re-fit the column widths now that the new fields/columns have content
$ws.Columns.Item(3).ColumnWidth = 5.833333333333333
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(6).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(8).ColumnWidth = 3.3333333333333335
$ws.Columns.Item(9).ColumnWidth = 10.166666666666666

# Update the selection to match the new active cell
$ws.Range("D15").Select()
